# Scheduled market-data refresh: updates currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) on the affected rows across all Leve sheets.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1611188.6
$ws.Range("I33").Value = 1764492.4
$ws.Range("J33").Value = 1500
$ws.Range("K33").Value = 1764492.4
$ws.Range("L33").Value = 1500
$ws.Range("M33").Value = -1764263.4
$ws.Range("N33").Value = -1958
$ws.Range("H41").Value = 685.625
$ws.Range("I41").Value = 520.75
$ws.Range("J41").Value = 850.5
$ws.Range("K41").Value = 520.75
$ws.Range("L41").Value = 850.5
$ws.Range("M41").Value = -80.75
$ws.Range("N41").Value = -1730.5
$ws.Range("H138").Value = 2990.2717
$ws.Range("J138").Value = 3126.2603
$ws.Range("L138").Value = 9378.7809
$ws.Range("N138").Value = -19658.7809

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 8061
$ws.Range("I19").Value = 5326.25
$ws.Range("K19").Value = 5326.25
$ws.Range("M19").Value = -5097.25
$ws.Range("H32").Value = 4824.328
$ws.Range("I32").Value = 3451.9023
$ws.Range("J32").Value = 7270.826
$ws.Range("K32").Value = 3451.9023
$ws.Range("L32").Value = 7270.826
$ws.Range("M32").Value = -3164.9023
$ws.Range("N32").Value = -7844.826
$ws.Range("H61").Value = 2279.4583
$ws.Range("I61").Value = 1672.5294
$ws.Range("J61").Value = 3753.4285
$ws.Range("K61").Value = 1672.5294
$ws.Range("L61").Value = 3753.4285
$ws.Range("M61").Value = -1460.5294
$ws.Range("N61").Value = -4177.4285
$ws.Range("H74").Value = 54091.402
$ws.Range("I74").Value = 3852.3572
$ws.Range("J74").Value = 265095.4
$ws.Range("K74").Value = 3852.3572
$ws.Range("L74").Value = 265095.4
$ws.Range("M74").Value = -2978.3572
$ws.Range("N74").Value = -266843.4
$ws.Range("H77").Value = 54091.402
$ws.Range("I77").Value = 3852.3572
$ws.Range("J77").Value = 265095.4
$ws.Range("K77").Value = 19261.786
$ws.Range("L77").Value = 1325477
$ws.Range("M77").Value = -14893.786
$ws.Range("N77").Value = -1334213
$ws.Range("H94").Value = 41142.57
$ws.Range("J94").Value = 41142.57
$ws.Range("L94").Value = 41142.57
$ws.Range("N94").Value = -42944.57
$ws.Range("H122").Value = 871616.9
$ws.Range("I122").Value = 3652.6667
$ws.Range("K122").Value = 10958.0001
$ws.Range("M122").Value = -8508.000100000001
$ws.Range("H132").Value = 2840
$ws.Range("I132").Value = 1700.2222
$ws.Range("J132").Value = 4122.25
$ws.Range("K132").Value = 5100.6666
$ws.Range("L132").Value = 12366.75
$ws.Range("M132").Value = -2570.6666
$ws.Range("N132").Value = -17426.75
$ws.Range("H136").Value = 2279.4583
$ws.Range("I136").Value = 1672.5294
$ws.Range("J136").Value = 3753.4285
$ws.Range("K136").Value = 5017.5882
$ws.Range("L136").Value = 11260.2855
$ws.Range("M136").Value = -2467.5882
$ws.Range("N136").Value = -16360.2855

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 50977.5
$ws.Range("J13").Value = 50977.5
$ws.Range("L13").Value = 50977.5
$ws.Range("N13").Value = -51313.5
$ws.Range("H86").Value = 3231673.2
$ws.Range("I86").Value = 4006558.5
$ws.Range("K86").Value = 4006558.5
$ws.Range("M86").Value = -4005435.5
$ws.Range("H89").Value = 3231673.2
$ws.Range("I89").Value = 4006558.5
$ws.Range("K89").Value = 20032792.5
$ws.Range("M89").Value = -20027176.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2715.7273
$ws.Range("I16").Value = 2325.5715
$ws.Range("K16").Value = 2325.5715
$ws.Range("M16").Value = -2038.5715
$ws.Range("H31").Value = 20550.812
$ws.Range("I31").Value = 1208.375
$ws.Range("J31").Value = 59235.688
$ws.Range("K31").Value = 1208.375
$ws.Range("L31").Value = 59235.688
$ws.Range("M31").Value = -913.375
$ws.Range("N31").Value = -59825.688
$ws.Range("H34").Value = 20550.812
$ws.Range("I34").Value = 1208.375
$ws.Range("J34").Value = 59235.688
$ws.Range("K34").Value = 1208.375
$ws.Range("L34").Value = 59235.688
$ws.Range("M34").Value = -1006.375
$ws.Range("N34").Value = -59639.688
$ws.Range("H97").Value = 28495.75
$ws.Range("J97").Value = 28495.75
$ws.Range("L97").Value = 28495.75
$ws.Range("N97").Value = -30477.75
$ws.Range("H109").Value = 34694.2
$ws.Range("J109").Value = 39617.75
$ws.Range("L109").Value = 39617.75
$ws.Range("N109").Value = -41697.75
$ws.Range("H113").Value = 2715.7273
$ws.Range("I113").Value = 2325.5715
$ws.Range("K113").Value = 2325.5715
$ws.Range("M113").Value = -155.5715

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1232.4615
$ws.Range("I7").Value = 2311.6
$ws.Range("J7").Value = 558
$ws.Range("K7").Value = 6934.799999999999
$ws.Range("L7").Value = 1674
$ws.Range("M7").Value = -6822.799999999999
$ws.Range("N7").Value = -1898
$ws.Range("H11").Value = 3540.8572
$ws.Range("I11").Value = 6945.6665
$ws.Range("J11").Value = 987.25
$ws.Range("K11").Value = 20836.9995
$ws.Range("L11").Value = 2961.75
$ws.Range("M11").Value = -20696.9995
$ws.Range("N11").Value = -3241.75
$ws.Range("H26").Value = 264.14285
$ws.Range("I26").Value = 331.4
$ws.Range("K26").Value = 994.1999999999999
$ws.Range("M26").Value = -706.1999999999999
$ws.Range("H57").Value = 7142.857
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H60").Value = 534.8
$ws.Range("I60").Value = 547
$ws.Range("K60").Value = 1641
$ws.Range("M60").Value = -1390
$ws.Range("H61").Value = 133.25
$ws.Range("I61").Value = 133.25
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 399.75
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -184.75
$ws.Range("N61").ClearContents()
$ws.Range("H87").Value = 12500
$ws.Range("I87").Value = 12500
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 37500
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -36252
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 12500
$ws.Range("I90").Value = 12500
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 112500
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -106260
$ws.Range("N90").ClearContents()
$ws.Range("H92").Value = 654.1111
$ws.Range("J92").Value = 669.8570999999999
$ws.Range("L92").Value = 2009.5713
$ws.Range("N92").Value = -4505.5713
$ws.Range("H107").Value = 194
$ws.Range("I107").Value = 175.95
$ws.Range("K107").Value = 527.8499999999999
$ws.Range("M107").Value = 1392.15
$ws.Range("H109").Value = 1592
$ws.Range("I109").Value = 1024
$ws.Range("J109").Value = 5000
$ws.Range("K109").Value = 3072
$ws.Range("L109").Value = 15000
$ws.Range("M109").Value = -2032
$ws.Range("N109").Value = -17080
$ws.Range("H136").Value = 2627
$ws.Range("I136").Value = 2627
$ws.Range("K136").Value = 7881
$ws.Range("M136").Value = -2781

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 8588583
$ws.Range("I126").Value = 4134720.2
$ws.Range("J126").Value = 20836706
$ws.Range("K126").Value = 12404160.6
$ws.Range("L126").Value = 62510118
$ws.Range("M126").Value = -12401690.6
$ws.Range("N126").Value = -62515058

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 24666.334
$ws.Range("I45").Value = 24666.334
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 24666.334
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -24259.334
$ws.Range("N45").ClearContents()
$ws.Range("H46").Value = 5354.0415
$ws.Range("I46").Value = 5036.364
$ws.Range("J46").Value = 5622.846
$ws.Range("K46").Value = 5036.364
$ws.Range("L46").Value = 5622.846
$ws.Range("M46").Value = -4848.364
$ws.Range("N46").Value = -5998.846
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 112323880
$ws.Range("I132").Value = 250000780
$ws.Range("J132").Value = 2182362.5
$ws.Range("K132").Value = 750002340
$ws.Range("L132").Value = 6547087.5
$ws.Range("M132").Value = -749999810
$ws.Range("N132").Value = -6552147.5

